$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 278.44
$ws.Range("I15").Value = 278.44
$ws.Range("K15").Value = 835.3199999999999
$ws.Range("M15").Value = -666.3199999999999
$ws.Range("H18").Value = 587.1
$ws.Range("I18").Value = 627.2857
$ws.Range("J18").Value = 493.33334
$ws.Range("K18").Value = 627.2857
$ws.Range("L18").Value = 493.33334
$ws.Range("M18").Value = -343.2857
$ws.Range("N18").Value = -1061.33334
$ws.Range("H19").Value = 14115.134
$ws.Range("I19").Value = 925.25
$ws.Range("J19").Value = 18911.455
$ws.Range("K19").Value = 925.25
$ws.Range("L19").Value = 18911.455
$ws.Range("M19").Value = -750.25
$ws.Range("N19").Value = -19261.455
$ws.Range("H64").Value = 2955.5557
$ws.Range("I64").Value = 2885.7144
$ws.Range("K64").Value = 2885.7144
$ws.Range("M64").Value = -2637.7144
$ws.Range("H67").Value = 2955.5557
$ws.Range("I67").Value = 2885.7144
$ws.Range("K67").Value = 2885.7144
$ws.Range("M67").Value = -2027.7144

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4519.6113
$ws.Range("I63").Value = 3400
$ws.Range("J63").Value = 4659.5625
$ws.Range("K63").Value = 3400
$ws.Range("L63").Value = 4659.5625
$ws.Range("M63").Value = -2714
$ws.Range("N63").Value = -6031.5625
$ws.Range("H66").Value = 4519.6113
$ws.Range("I66").Value = 3400
$ws.Range("J66").Value = 4659.5625
$ws.Range("K66").Value = 17000
$ws.Range("L66").Value = 23297.8125
$ws.Range("M66").Value = -13568
$ws.Range("N66").Value = -30161.8125
$ws.Range("H132").Value = 3116.5344
$ws.Range("I132").Value = 2703.0557
$ws.Range("J132").Value = 3793.1365
$ws.Range("K132").Value = 8109.1671
$ws.Range("L132").Value = 11379.4095
$ws.Range("M132").Value = -5579.1671
$ws.Range("N132").Value = -16439.4095

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 43480080
$ws.Range("I20").Value = 1878
$ws.Range("J20").Value = 100001740
$ws.Range("K20").Value = 1878
$ws.Range("L20").Value = 100001740
$ws.Range("M20").Value = -1631
$ws.Range("N20").Value = -100002234
$ws.Range("H99").Value = 1186.4117
$ws.Range("I99").Value = 1046.5834
$ws.Range("J99").Value = 1522
$ws.Range("K99").Value = 1046.5834
$ws.Range("L99").Value = 1522
$ws.Range("M99").Value = 451.4166
$ws.Range("N99").Value = -4518
$ws.Range("H100").Value = 67661
$ws.Range("J100").Value = 67661
$ws.Range("L100").Value = 67661
$ws.Range("N100").Value = -69825
$ws.Range("H134").Value = 1895.8445
$ws.Range("I134").Value = 1674.6786
$ws.Range("J134").Value = 2260.1177
$ws.Range("K134").Value = 5024.0358
$ws.Range("L134").Value = 6780.353099999999
$ws.Range("M134").Value = -2489.0358
$ws.Range("N134").Value = -11850.3531

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()
$ws.Range("H62").Value = 4365.467
$ws.Range("I62").Value = 4488.4307
$ws.Range("J62").Value = 3566.2
$ws.Range("K62").Value = 4488.4307
$ws.Range("L62").Value = 3566.2
$ws.Range("M62").Value = -3864.4307
$ws.Range("N62").Value = -4814.2
$ws.Range("H65").Value = 4365.467
$ws.Range("I65").Value = 4488.4307
$ws.Range("J65").Value = 3566.2
$ws.Range("K65").Value = 22442.1535
$ws.Range("L65").Value = 17831
$ws.Range("M65").Value = -19322.1535
$ws.Range("N65").Value = -24071
$ws.Range("H99").Value = 2457.0454
$ws.Range("I99").Value = 2299.5833
$ws.Range("J99").Value = 2516.0938
$ws.Range("K99").Value = 2299.5833
$ws.Range("L99").Value = 2516.0938
$ws.Range("M99").Value = -801.5832999999998
$ws.Range("N99").Value = -5512.093800000001
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()
$ws.Range("H126").Value = 2457.0454
$ws.Range("I126").Value = 2299.5833
$ws.Range("J126").Value = 2516.0938
$ws.Range("K126").Value = 6898.749899999999
$ws.Range("L126").Value = 7548.2814
$ws.Range("M126").Value = -4428.749899999999
$ws.Range("N126").Value = -12488.2814
$ws.Range("H127").Value = 88780
$ws.Range("J127").Value = 88780
$ws.Range("L127").Value = 88780
$ws.Range("N127").Value = -98700
$ws.Range("H134").Value = 19248.5
$ws.Range("I134").Value = 27372.75
$ws.Range("K134").Value = 82118.25
$ws.Range("M134").Value = -79583.25

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 65.791664
$ws.Range("I12").Value = 35.3125
$ws.Range("J12").Value = 126.75
$ws.Range("K12").Value = 105.9375
$ws.Range("L12").Value = 380.25
$ws.Range("M12").Value = 67.0625
$ws.Range("N12").Value = -726.25
$ws.Range("H64").Value = 1335469.2
$ws.Range("I64").Value = 1180
$ws.Range("J64").Value = 1669041.6
$ws.Range("K64").Value = 3540
$ws.Range("L64").Value = 5007124.800000001
$ws.Range("M64").Value = -3270
$ws.Range("N64").Value = -5007664.800000001
$ws.Range("H67").Value = 1335469.2
$ws.Range("I67").Value = 1180
$ws.Range("J67").Value = 1669041.6
$ws.Range("K67").Value = 3540
$ws.Range("L67").Value = 5007124.800000001
$ws.Range("M67").Value = -2604
$ws.Range("N67").Value = -5008996.800000001
$ws.Range("H75").Value = 999
$ws.Range("I75").Value = 999
$ws.Range("K75").Value = 2997
$ws.Range("M75").Value = -1999
$ws.Range("H78").Value = 999
$ws.Range("I78").Value = 999
$ws.Range("K78").Value = 8991
$ws.Range("M78").Value = -3999
$ws.Range("H137").Value = 47902.293
$ws.Range("J137").Value = 128637.5
$ws.Range("L137").Value = 385912.5
$ws.Range("N137").Value = -396112.5
$ws.Range("H140").Value = 1757.2646
$ws.Range("I140").Value = 1150.35
$ws.Range("K140").Value = 3451.05
$ws.Range("M140").Value = 1728.95

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 5402.75
$ws.Range("I70").Value = 5336.095
$ws.Range("J70").Value = 5602.7144
$ws.Range("K70").Value = 5336.095
$ws.Range("L70").Value = 5602.7144
$ws.Range("M70").Value = -5066.095
$ws.Range("N70").Value = -6142.7144
$ws.Range("H73").Value = 5402.75
$ws.Range("I73").Value = 5336.095
$ws.Range("J73").Value = 5602.7144
$ws.Range("K73").Value = 5336.095
$ws.Range("L73").Value = 5602.7144
$ws.Range("M73").Value = -4400.095
$ws.Range("N73").Value = -7474.7144
$ws.Range("H132").Value = 25002292
$ws.Range("I132").Value = 30304808
$ws.Range("J132").Value = 4717.7144
$ws.Range("K132").Value = 90914424
$ws.Range("L132").Value = 14153.1432
$ws.Range("M132").Value = -90911894
$ws.Range("N132").Value = -19213.1432

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2818.3333
$ws.Range("I40").Value = 3227.5
$ws.Range("K40").Value = 3227.5
$ws.Range("M40").Value = -3091.5
$ws.Range("H140").Value = 57062.75
$ws.Range("J140").Value = 57062.75
$ws.Range("L140").Value = 57062.75
$ws.Range("N140").Value = -67422.75

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws.Range("H107").Value = 1408.3334
$ws.Range("I107").Value = 1237.5
$ws.Range("J107").Value = 1750
$ws.Range("K107").Value = 3712.5
$ws.Range("L107").Value = 5250
$ws.Range("M107").Value = -1792.5
$ws.Range("N107").Value = -9090
$ws.Range("H132").Value = 4067846.2
$ws.Range("I132").Value = 2621.7932
$ws.Range("J132").Value = 13892139
$ws.Range("K132").Value = 7865.3796
$ws.Range("L132").Value = 41676417
$ws.Range("M132").Value = -5335.3796
$ws.Range("N132").Value = -41681477
